$wb = $excel.ActiveWorkbook

# Map of row -> new F-column value, shared by both "展览" and "全部类型" sheets
$updates = @{
    3  = 384
    4  = 1485
    5  = 8718
    6  = 88
    11 = 14
    12 = 3567
    13 = 49
    15 = 77
    16 = 1164
    17 = 146
    18 = 1115
    20 = 202
    21 = 2324
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# Last row (row 22 on "展览", row 23 on "全部类型") also changes 53 -> 54
$wb.Worksheets.Item("展览").Range("F22").Value = 54
$wb.Worksheets.Item("全部类型").Range("F23").Value = 54
